# Update newly added iAuthor TC's credentials on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Candidate ID 231011253 -> 231027165)
$ws.Range("A2").Value = "PGBfK762"
$ws.Range("B2").Value = 231027165
$ws.Range("C2").Value = "fokfqah75"
$ws.Range("D2").Value = "uP!N9y#4"
$ws.Range("F2").Value = "rqilgaUh"
$ws.Range("G2").Value = "ELjT"

# Row 3 (Candidate ID 231011252 -> 231027164)
$ws.Range("A3").Value = "ZKoJP429"
$ws.Range("B3").Value = 231027164
$ws.Range("C3").Value = "xkxhpzg51"
$ws.Range("D3").Value = "E%Yse$74"
$ws.Range("F3").Value = "xKZdAMDi"
$ws.Range("G3").Value = "axSX"
